$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()

$ws.Range("C4").Value = -0.8792832172735632
$ws.Range("E4").Value = -0.1037449741509211

$ws.Range("C5").Value = 0.9337833426867226

$ws.Range("C6").Value = 2.791140000794257
$ws.Range("E6").Value = 1.722110645261954

$ws.Range("C7").Value = 0.4451370000809529
$ws.Range("E7").Value = 0.6480763427742176

$ws.Range("E8").Value = 1.223618887196509

$ws.Range("E10").Value = 1.582979977679555

$ws.Range("C11").Value = 2.2044495746113
$ws.Range("E11").Value = 1.401113624217065

$ws.Range("E12").Value = 2.357704431248386

$ws.Range("E13").Value = 2.36261304543155

$ws.Range("E14").Value = 0.7756897792100093

$ws.Range("C15").Value = -3.258619210312896
$ws.Range("E15").Value = -1.049961713694159

$ws.Range("C16").Value = 0.4255262881966759
$ws.Range("E16").Value = 2.734996705911397

$ws.Range("E17").Value = 1.813346177122321

$ws.Range("C18").Value = -0.2814561130375703

$ws.Range("C19").Value = -0.6470065423293758
$ws.Range("E19").Value = 0.9227184786156251
